{"js": "// Remove the empty \"ACADEMIC PROJECTS / None\", \"HONORS AND AWARDS / None\",\n// and \"ACTIVITIES / None\" placeholder sections from the resume (they had\n// no real content), along with their surrounding blank-line paragraphs.\n//\n// Resulting paragraph flow:\n//   ... Cumulative GPA: 4.0 | (blank) | WORK EXPERIENCE | Strategic Staffing...\n//   ... Successfully handled...S3. | (blank) | SKILLS AND CERTIFICATES ...\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst count = items.length;\n\n// Collect the (0-based) indices of every paragraph we need to delete, then\n// delete from the bottom of the document upward so earlier indices stay valid.\nconst toDelete = [];\n\nfor (let i = 0; i < count; i++) {\n  const text = items[i].text.trim();\n\n  if (text === \"ACADEMIC PROJECTS\") {\n    // Heading + its \"None\" line + the blank line that follows it.\n    toDelete.push(i);\n    if (i + 1 < count && items[i + 1].text.trim() === \"None\") {\n      toDelete.push(i + 1);\n    }\n    if (i + 2 < count && items[i + 2].text.trim() === \"\") {\n      toDelete.push(i + 2);\n    }\n  } else if (text === \"HONORS AND AWARDS\") {\n    // The blank line before the heading, the heading itself, its \"None\"\n    // line, and the blank line that follows it.\n    if (i - 1 >= 0 && items[i - 1].text.trim() === \"\") {\n      toDelete.push(i - 1);\n    }\n    toDelete.push(i);\n    if (i + 1 < count && items[i + 1].text.trim() === \"None\") {\n      toDelete.push(i + 1);\n    }\n    if (i + 2 < count && items[i + 2].text.trim() === \"\") {\n      toDelete.push(i + 2);\n    }\n  } else if (text === \"ACTIVITIES\") {\n    // Heading itself and its trailing \"None\" line.\n    toDelete.push(i);\n    if (i + 1 < count && items[i + 1].text.trim() === \"None\") {\n      toDelete.push(i + 1);\n    }\n  }\n}\n\nconst uniqueSorted = Array.from(new Set(toDelete)).sort((a, b) => b - a);\nfor (const idx of uniqueSorted) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty \"ACADEMIC PROJECTS / None\", \"HONORS AND AWARDS / None\",\n# and \"ACTIVITIES / None\" placeholder sections from the resume (they had\n# no real content), along with their surrounding blank-line paragraphs.\n#\n# Resulting paragraph flow:\n#   ... Cumulative GPA: 4.0 | (blank) | WORK EXPERIENCE | Strategic Staffing...\n#   ... Successfully handled...S3. | (blank) | SKILLS AND CERTIFICATES ...\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]13, [char]7).Trim()\n}\n\n$count = $d.Paragraphs.Count\n\n# Snapshot paragraph texts first (indices are 1-based, like real Word COM).\n$texts = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $texts += ,(Get-ParaText $d.Paragraphs($i))\n}\n\n$toDelete = @()\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $texts[$i - 1]\n\n    if ($text -eq \"ACADEMIC PROJECTS\") {\n        # Heading + its \"None\" line + the blank line that follows it.\n        $toDelete += $i\n        if ($i + 1 -le $count -and $texts[$i] -eq \"None\") {\n            $toDelete += ($i + 1)\n        }\n        if ($i + 2 -le $count -and $texts[$i + 1] -eq \"\") {\n            $toDelete += ($i + 2)\n        }\n    } elseif ($text -eq \"HONORS AND AWARDS\") {\n        # The blank line before the heading, the heading itself, its \"None\"\n        # line, and the blank line that follows it.\n        if ($i - 1 -ge 1 -and $texts[$i - 2] -eq \"\") {\n            $toDelete += ($i - 1)\n        }\n        $toDelete += $i\n        if ($i + 1 -le $count -and $texts[$i] -eq \"None\") {\n            $toDelete += ($i + 1)\n        }\n        if ($i + 2 -le $count -and $texts[$i + 1] -eq \"\") {\n            $toDelete += ($i + 2)\n        }\n    } elseif ($text -eq \"ACTIVITIES\") {\n        # Heading itself and its trailing \"None\" line.\n        $toDelete += $i\n        if ($i + 1 -le $count -and $texts[$i] -eq \"None\") {\n            $toDelete += ($i + 1)\n        }\n    }\n}\n\n$uniqueSorted = $toDelete | Select-Object -Unique | Sort-Object -Descending\n\nforeach ($idx in $uniqueSorted) {\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
